$d = $word.ActiveDocument

# Locate the paragraph that ends with "...já que o prazo já é curto." (the
# SCRUM paragraph) so we can insert the new content right after it,
# regardless of its exact paragraph index.
$searchRange = $d.Content
$found = $searchRange.Find.Execute(
    "já que o prazo já é curto.",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Anchor paragraph not found"
}

$anchorIndex = $searchRange.Paragraphs(1).Index

# Split off a brand-new (empty) paragraph right after the anchor paragraph
# first. Running InsertXML directly on a range collapsed at the anchor
# paragraph's own end would overwrite the anchor paragraph itself, so we
# grow the document by one empty paragraph and target the fresh one.
$anchorPara = $d.Paragraphs($anchorIndex)
$tail = $anchorPara.Range.Duplicate
$tail.Collapse(0)
$tail.InsertParagraphAfter()

$newPara = $d.Paragraphs($anchorIndex + 1)
$insertionPoint = $newPara.Range.Duplicate
$insertionPoint.Collapse(1)

# Build the new paragraphs as a WordprocessingML fragment so formatting
# details (proofing marks, preserved whitespace, etc.) match exactly.
$fragment = @'
<w:p/><w:p><w:r><w:t xml:space="preserve">Eu organizaria a equipe com o modelo ágil </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>( SCRUM</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> ), assim teria 3 categorias :</w:t></w:r></w:p><w:p><w:proofErr w:type="gramStart"/><w:r><w:t>PO :</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> também chamado de </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Product</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Owner</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, ele é o responsável por garantir retorno de investimento. Ele é quem sabe </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>oque</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> o cliente quer, existe apenas um por projeto.</w:t></w:r></w:p><w:p><w:proofErr w:type="gramStart"/><w:r><w:t>SM :</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> também chamado de Scrum Master, ele é o responsável por facilitar a vida do time, tirando da frente problemas desnecessários, além disso ele garante que o Scrum flua como deve e proteger o time de interferências externas.</w:t></w:r></w:p><w:p><w:r><w:t>Time: Produzir os produtos, além de definir metas e ter autonomia no gerenciamento.</w:t></w:r></w:p>
'@

$insertionPoint.InsertXML($fragment)
